$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '24.431.11'
$ws.Cells.Item(2, 5).Value = '  -1.64%  '
$ws.Cells.Item(3, 4).Value = '1.686.43'
$ws.Cells.Item(3, 5).Value = '  -1.21%  '
$ws.Cells.Item(4, 5).Value = '  -0.16%  '
$ws.Cells.Item(5, 4).Value = "'315.94"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  -0.15%  '
$ws.Cells.Item(6, 4).Value = "'0.9999"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = '  -0.12%  '
$ws.Cells.Item(7, 4).Value = "'0.3909"
$ws.Cells.Item(7, 4).Style = "Normal"
$ws.Cells.Item(7, 5).Value = '  -0.76%  '
$ws.Cells.Item(8, 4).Value = "'0.4029"
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Cells.Item(8, 5).Value = '  -0.65%  '
$ws.Cells.Item(9, 4).Value = "'1.486"
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Cells.Item(9, 5).Value = '  -2.12%  '
$ws.Cells.Item(10, 4).Value = "'0.9997"
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).Value = '  -0.21%  '
$ws.Cells.Item(11, 4).Value = "'52.48"
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 5).Value = '  -2.13%  '
$ws.Cells.Item(12, 4).Value = "'0.08771"
$ws.Cells.Item(12, 4).Style = "Normal"
$ws.Cells.Item(12, 5).Value = '  -1.57%  '
$ws.Cells.Item(13, 4).Value = "'26.52"
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Cells.Item(13, 5).Value = '  +12.32%  '
$ws.Cells.Item(14, 4).Value = "'7.443"
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).Value = '  +1.69%  '
$ws.Cells.Item(15, 4).Value = "'8.155"
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Cells.Item(15, 5).Value = '  +1.43%  '
$ws.Cells.Item(16, 4).Value = "'0.00001344"
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Cells.Item(16, 5).Value = '  +0.92%  '
$ws.Cells.Item(17, 4).Value = '1.686.99'
$ws.Cells.Item(17, 5).Value = '  -0.96%  '
$ws.Cells.Item(18, 4).Value = "'98.04"
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Cells.Item(18, 5).Value = '  -2.36%  '
$ws.Cells.Item(19, 4).Value = "'0.07234"
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Cells.Item(19, 5).Value = '  +2.64%  '
$ws.Cells.Item(20, 4).Value = "'20.29"
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value = '  +2.86%  '
$ws.Cells.Item(21, 4).Value = "'7.262"
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).Value = '  +2.49%  '
$ws.Cells.Item(22, 4).Value = "'0.9996"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = '  -0.11%  '
$ws.Cells.Item(23, 4).Value = "'14.27"
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).Value = '  -2.07%  '
$ws.Cells.Item(24, 4).Value = '24.419.51'
$ws.Cells.Item(24, 5).Value = '  -1.66%  '
$ws.Cells.Item(25, 4).Value = "'3.037"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = '  -5.57%  '
$ws.Cells.Item(26, 4).Value = "'2.336"
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(26, 5).Value = '  -1.11%  '
$ws.Cells.Item(27, 4).Value = "'22.57"
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(27, 5).Value = '  -1.40%  '
$ws.Cells.Item(28, 4).Value = "'167.13"
$ws.Cells.Item(28, 4).Style = "Normal"
$ws.Cells.Item(28, 5).Value = '  +3.22%  '
$ws.Cells.Item(29, 4).Value = "'8.462"
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Cells.Item(29, 5).Value = '  +0.73%  '
$ws.Cells.Item(30, 4).Value = "'5.335"
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Cells.Item(30, 5).Value = '  +3.07%  '
$ws.Cells.Item(31, 4).Value = "'138.11"
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Cells.Item(31, 5).Value = '  +0.59%  '
$ws.Cells.Item(32, 4).Value = '1.867.26'
$ws.Cells.Item(32, 5).Value = '  -1.09%  '
$ws.Cells.Item(33, 4).Value = "'0.08769"
$ws.Cells.Item(33, 4).Style = "Normal"
$ws.Cells.Item(33, 5).Value = '  -1.54%  '
$ws.Cells.Item(34, 4).Value = "'7.271"
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Cells.Item(34, 5).Value = '  -3.42%  '
$ws.Cells.Item(35, 2).Value = 'ImmutableX'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Cells.Item(35, 4).Value = "'1.047"
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Cells.Item(35, 5).Value = '  -3.71%  '
$ws.Cells.Item(36, 2).Value = 'WEMIXTOKEN'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Cells.Item(36, 4).Value = "'2.093"
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Cells.Item(36, 5).Value = '  +5.79%  '
$ws.Cells.Item(37, 4).Value = "'0.03008"
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(37, 5).Value = '  +8.85%  '
$ws.Cells.Item(38, 4).Value = "'0.2782"
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Cells.Item(38, 5).Value = '  +1.01%  '
$ws.Cells.Item(39, 4).Value = "'10.82"
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Cells.Item(39, 5).Value = '  -3.21%  '
$ws.Cells.Item(40, 4).Value = "'0.09138"
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).Value = '  -1.08%  '
$ws.Cells.Item(41, 4).Value = "'14.16"
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).Value = '  -2.42%  '
$ws.Cells.Item(42, 4).Value = "'0.8023"
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value = '  +3.77%  '
$ws.Cells.Item(43, 4).Value = "'1.472"
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).Value = '  +0.68%  '
$ws.Cells.Item(44, 4).Value = "'17.60"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).Value = '  +10.97%  '
$ws.Cells.Item(45, 4).Value = "'2.653"
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value = '  +2.97%  '
$ws.Cells.Item(46, 4).Value = "'0.7234"
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Value = '  +0.24%  '
$ws.Cells.Item(47, 4).Value = "'4.259"
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = '  +1.14%  '
$ws.Cells.Item(48, 4).Value = "'1.414"
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(48, 5).Value = '  +6.82%  '
$ws.Cells.Item(49, 4).Value = "'0.9994"
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value = '  -0.18%  '
$ws.Cells.Item(50, 4).Value = "'139.08"
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = '  -1.15%  '
$ws.Cells.Item(51, 4).Value = "'0.08080"
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value = '  +0.80%  '
